# Generate Report for Handoff
#
# The 5afef0f6-... handback file turned out to be stale / superseded, so its
# status flips from "Handed back: in sync with en-US" to "Ready for handoff"
# (with refreshed handoff timestamps and a new "version not latest" error
# message), while the d992f641-... entry keeps its original status. Rows in
# each sheet are kept sorted by file name, so the two data rows swap places
# (the d992f641 entry moves up to row 2, 5afef0f6 moves down to row 3).

$wb = $excel.ActiveWorkbook

$fileA = "5afef0f6-2833-4e05-be18-cb778151c15b.md"
$fileB = "d992f641-649d-4f82-9446-27d81d55f3e9.md"

# ---------------------------------------------------------------------------
# Sheet "Overview": row 2 <-> row 3 identity swap, plus refreshed status for
# the 5afef0f6 row (now row 3).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $fileB
$wsOverview.Range("B2").Value = "e2e\" + $fileB
$wsOverview.Range("A3").Value = $fileA
$wsOverview.Range("B3").Value = "e2e\" + $fileA

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-15 10:46:58"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Row -eq 2) {
        $hl.TextToDisplay = "e2e\" + $fileB
    } elseif ($hl.Range.Row -eq 3) {
        $hl.TextToDisplay = "e2e\" + $fileA
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $fileB
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("G2").Value = "d992f641-649d-4f82-9446-27d81d55f3e9.2a45df4773e7d5df37c9e93e673b7802772a01a5.zh-cn.xlf"
$wsZh.Range("I2").Value = $fileB
$wsZh.Range("J2").Value = "d992f641-649d-4f82-9446-27d81d55f3e9.2a45df4773e7d5df37c9e93e673b7802772a01a5.zh-cn.xlf"

$wsZh.Range("A3").Value = $fileA
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("G3").Value = "5afef0f6-2833-4e05-be18-cb778151c15b.86715cbe9b9cc0be8b71dcb30afc0f3a0eb363b4.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-15 10:46:54"
$wsZh.Range("I3").Value = $fileA
$wsZh.Range("J3").Value = "5afef0f6-2833-4e05-be18-cb778151c15b.86715cbe9b9cc0be8b71dcb30afc0f3a0eb363b4.zh-cn.xlf"
$wsZh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e20e5adf5d78b823cb5ebf1531a13acb3349efd/e2e/5afef0f6-2833-4e05-be18-cb778151c15b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f5c58d9baada591e9f7174f254a6f0aeccbefe9/e2e/5afef0f6-2833-4e05-be18-cb778151c15b.md."

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Row -eq 2) {
        $hl.TextToDisplay = $fileB
    } elseif ($hl.Range.Row -eq 3) {
        $hl.TextToDisplay = $fileA
    }
}

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $fileB
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("G2").Value = "d992f641-649d-4f82-9446-27d81d55f3e9.2a45df4773e7d5df37c9e93e673b7802772a01a5.de-de.xlf"
$wsDe.Range("I2").Value = $fileB
$wsDe.Range("J2").Value = "d992f641-649d-4f82-9446-27d81d55f3e9.2a45df4773e7d5df37c9e93e673b7802772a01a5.de-de.xlf"

$wsDe.Range("A3").Value = $fileA
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("G3").Value = "5afef0f6-2833-4e05-be18-cb778151c15b.86715cbe9b9cc0be8b71dcb30afc0f3a0eb363b4.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-15 10:46:58"
$wsDe.Range("I3").Value = $fileA
$wsDe.Range("J3").Value = "5afef0f6-2833-4e05-be18-cb778151c15b.86715cbe9b9cc0be8b71dcb30afc0f3a0eb363b4.de-de.xlf"
$wsDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e20e5adf5d78b823cb5ebf1531a13acb3349efd/e2e/5afef0f6-2833-4e05-be18-cb778151c15b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f5c58d9baada591e9f7174f254a6f0aeccbefe9/e2e/5afef0f6-2833-4e05-be18-cb778151c15b.md."

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Row -eq 2) {
        $hl.TextToDisplay = $fileB
    } elseif ($hl.Range.Row -eq 3) {
        $hl.TextToDisplay = $fileA
    }
}

$wsDe.Columns.Item(16).ColumnWidth = 39.17
